# "Generate Report for Handback" — refresh the handback status/report sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: Overview ---
$ovw = $wb.Worksheets.Item("Overview")
# 1f748b26-....md is no longer in sync with en-US
$ovw.Range("E2").Value = "Handed back: not in sync with en-US"
$ovw.Range("F2").Value = "Handed back: not in sync with en-US"
# 26dc01cb-....md got a fresh handback report timestamp
$ovw.Range("G3").Value = "2017-02-22 08:19:30"
$ovw.Columns.Item(5).ColumnWidth = 32.666666666666664
$ovw.Columns.Item(6).ColumnWidth = 32.666666666666664

# --- Sheet 2: zh-cn ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "Handed back: not in sync with en-US"
$zh.Range("H3").Value = "2017-02-22 08:19:13"
$zh.Range("L3").Value = "2017-02-22 08:20:14"
$zh.Columns.Item(3).ColumnWidth = 32.666666666666664

# --- Sheet 3: de-de ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "Handed back: not in sync with en-US"
$de.Range("H3").Value = "2017-02-22 08:19:30"
$de.Range("L3").Value = "2017-02-22 08:20:35"
$de.Columns.Item(3).ColumnWidth = 32.666666666666664
